# Update the "Caso" (case) id for the "Pedraza Manuela 4101" record from -520 to 6538
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A59").Value = "6538"

# Remove the rows that were pulled out of the dataset.
# Deleted from the bottom up so earlier row numbers stay valid.
$ws.Rows.Item(52).EntireRow.Delete()   # Caso 6336 - PARAGUAY 4291
$ws.Rows.Item(50).EntireRow.Delete()   # Caso 6331 - PARAGUAY 4259
$ws.Rows.Item(28).EntireRow.Delete()   # Caso 804922192 - Paraguay 4657
